# Insert a new data row at row 224 (pushing existing rows 224-271 down to
# 225-272) on the single worksheet, then populate the new row with its
# values. This reproduces the diff, where a weekly price-report row was
# inserted in the middle of the "Betarraga" series and the dimension grew
# from A1:R271 to A1:R272.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 224..271 down to 225..272, inheriting row 224's existing
# formatting (e.g. the date style on column D) for the newly opened row.
$ws.Rows.Item(224).Insert()

# Populate the newly inserted row 224 with the new record.
$ws.Cells.Item(224, 1).Value2 = 4
$ws.Cells.Item(224, 2).Value2 = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(224, 3).Value2 = 'Los Lagos'
$ws.Cells.Item(224, 4).Value2 = 44641
$ws.Cells.Item(224, 5).Value2 = 10
$ws.Cells.Item(224, 6).Value2 = 100114014
$ws.Cells.Item(224, 7).Value2 = 'Betarraga'
$ws.Cells.Item(224, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(224, 9).Value2 = 'Primera'
$ws.Cells.Item(224, 10).Value2 = 500
$ws.Cells.Item(224, 11).Value2 = 1000
$ws.Cells.Item(224, 12).Value2 = 1000
$ws.Cells.Item(224, 13).Value2 = 1000
$ws.Cells.Item(224, 14).Value2 = '$/paquete 5 unidades'
$ws.Cells.Item(224, 15).Value2 = 'Región del Maule'
$ws.Cells.Item(224, 16).Value2 = 200
$ws.Cells.Item(224, 17).Value2 = 5
$ws.Cells.Item(224, 18).Value2 = 'Hortaliza'
